# Add a second slide ("Create folder and see it in git window") to the
# presentation, using the "Title and Content" layout (same layout family
# used by slideLayout2.xml -> ppLayoutText / obj type).
#
# This corresponds to the OOXML diff that:
#   - appends <p:sldId id="257" .../> to the presentation's sldIdLst
#   - adds a new ppt/slides/slide2.xml with a Title placeholder and a
#     Content placeholder, each populated with text.

$p = $ppt.ActivePresentation

# Index 2 == insert after the existing (first) slide.
# Layout 2 == ppLayoutText ("Title and Content"), the same family of layout
# used by the new slide in the target deck (Title placeholder + single
# unnumbered content/body placeholder).
$s = $p.Slides.Add(2, 2)

# --- Title placeholder ------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Create folder and see it in "
[void]$title.InsertAfter("git")
[void]$title.InsertAfter(" window")

# Shrink text to fit the placeholder, matching the normAutofit seen on the
# title shape in the target slide.
$s.Shapes.Item(1).TextFrame.AutoSize = 2

# --- Content / body placeholder ----------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "We can create any new folder to our any drive and see it on "
[void]$body.InsertAfter("git")
[void]$body.InsertAfter(" command prompt.")
[void]$body.InsertAfter("`rAs well as add that folder to our repository")
[void]$body.InsertAfter("`rUsing some commands are as follow")
